$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1250.6154
$ws.Range("J19").Value = 1889.75
$ws.Range("L19").Value = 1889.75
$ws.Range("N19").Value = -2239.75

$ws.Range("H38").Value = 1590.9412
$ws.Range("I38").Value = 147.55556
$ws.Range("J38").Value = 3214.75
$ws.Range("K38").Value = 442.66668
$ws.Range("L38").Value = 9644.25
$ws.Range("M38").Value = -70.66668000000004
$ws.Range("N38").Value = -10388.25

$ws.Range("H40").Value = 1994.2307
$ws.Range("I40").Value = 2078.0625
$ws.Range("J40").Value = 1860.1
$ws.Range("K40").Value = 2078.0625
$ws.Range("L40").Value = 1860.1
$ws.Range("M40").Value = -1903.0625
$ws.Range("N40").Value = -2210.1

$ws.Range("H43").Value = 950
$ws.Range("I43").Value = 900
$ws.Range("J43").Value = 1000
$ws.Range("K43").Value = 900
$ws.Range("L43").Value = 1000
$ws.Range("M43").Value = -831
$ws.Range("N43").Value = -1138

$ws.Range("H92").Value = 483687.1
$ws.Range("I92").Value = 654086.9399999999
$ws.Range("J92").Value = 887.5
$ws.Range("K92").Value = 654086.9399999999
$ws.Range("L92").Value = 887.5
$ws.Range("M92").Value = -652838.9399999999
$ws.Range("N92").Value = -3383.5

$ws.Range("H129").Value = 1451.5405
$ws.Range("I129").Value = 506.33334
$ws.Range("J129").Value = 1634.4839
$ws.Range("K129").Value = 1519.00002
$ws.Range("L129").Value = 4903.4517
$ws.Range("M129").Value = 3480.99998
$ws.Range("N129").Value = -14903.4517

$ws.Range("H132").Value = 265603.78
$ws.Range("I132").Value = 290006.34
$ws.Range("J132").Value = 60622.4
$ws.Range("K132").Value = 870019.02
$ws.Range("L132").Value = 181867.2
$ws.Range("M132").Value = -867489.02
$ws.Range("N132").Value = -186927.2

$ws.Range("H135").Value = 1500.7142
$ws.Range("I135").Value = 1287.5
$ws.Range("K135").Value = 11587.5
$ws.Range("M135").Value = -9052.5

$ws.Range("H137").Value = 19231930
$ws.Range("I137").Value = 35715260
$ws.Range("J137").Value = 1380.0416
$ws.Range("K137").Value = 107145780
$ws.Range("L137").Value = 4140.1248
$ws.Range("M137").Value = -107143230
$ws.Range("N137").Value = -9240.1248

$ws.Range("H138").Value = 3790142.5
$ws.Range("I138").Value = 732534.75
$ws.Range("J138").Value = 7939753
$ws.Range("K138").Value = 2197604.25
$ws.Range("L138").Value = 23819259
$ws.Range("M138").Value = -2192464.25
$ws.Range("N138").Value = -23829539

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19325.416
$ws.Range("I32").Value = 2587.9124
$ws.Range("K32").Value = 2587.9124
$ws.Range("M32").Value = -2300.9124

$ws.Range("H61").Value = 2051.628
$ws.Range("I61").Value = 1367.9429
$ws.Range("J61").Value = 5042.75
$ws.Range("K61").Value = 1367.9429
$ws.Range("L61").Value = 5042.75
$ws.Range("M61").Value = -1155.9429
$ws.Range("N61").Value = -5466.75

$ws.Range("H74").Value = 4167.628
$ws.Range("I74").Value = 1135.5758
$ws.Range("J74").Value = 14173.4
$ws.Range("K74").Value = 1135.5758
$ws.Range("L74").Value = 14173.4
$ws.Range("M74").Value = -261.5758000000001
$ws.Range("N74").Value = -15921.4

$ws.Range("H77").Value = 4167.628
$ws.Range("I77").Value = 1135.5758
$ws.Range("J77").Value = 14173.4
$ws.Range("K77").Value = 5677.879000000001
$ws.Range("L77").Value = 70867
$ws.Range("M77").Value = -1309.879000000001
$ws.Range("N77").Value = -79603

$ws.Range("H123").Value = 43624.625
$ws.Range("J123").Value = 43624.625
$ws.Range("L123").Value = 43624.625
$ws.Range("N123").Value = -53424.625

$ws.Range("H132").Value = 3024.4285
$ws.Range("I132").Value = 2478.6086
$ws.Range("J132").Value = 4070.5833
$ws.Range("K132").Value = 7435.825800000001
$ws.Range("L132").Value = 12211.7499
$ws.Range("M132").Value = -4905.825800000001
$ws.Range("N132").Value = -17271.7499

$ws.Range("H136").Value = 2051.628
$ws.Range("I136").Value = 1367.9429
$ws.Range("J136").Value = 5042.75
$ws.Range("K136").Value = 4103.8287
$ws.Range("L136").Value = 15128.25
$ws.Range("M136").Value = -1553.8287
$ws.Range("N136").Value = -20228.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3356.8975
$ws.Range("I134").Value = 2332.28
$ws.Range("K134").Value = 6996.84
$ws.Range("M134").Value = -4461.84

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1227.717
$ws.Range("I58").Value = 656.2105
$ws.Range("J58").Value = 2675.5334
$ws.Range("K58").Value = 656.2105
$ws.Range("L58").Value = 2675.5334
$ws.Range("M58").Value = -453.2105
$ws.Range("N58").Value = -3081.5334

$ws.Range("H132").Value = 1977.3492
$ws.Range("I132").Value = 1506.66
$ws.Range("J132").Value = 3787.6924
$ws.Range("K132").Value = 4519.98
$ws.Range("L132").Value = 11363.0772
$ws.Range("M132").Value = -1989.98
$ws.Range("N132").Value = -16423.0772

$ws.Range("H134").Value = 2100.492
$ws.Range("I134").Value = 1287.46
$ws.Range("J134").Value = 5227.5386
$ws.Range("K134").Value = 3862.38
$ws.Range("L134").Value = 15682.6158
$ws.Range("M134").Value = -1327.38
$ws.Range("N134").Value = -20752.6158

$ws.Range("H136").Value = 1227.717
$ws.Range("I136").Value = 656.2105
$ws.Range("J136").Value = 2675.5334
$ws.Range("K136").Value = 1968.6315
$ws.Range("L136").Value = 8026.600199999999
$ws.Range("M136").Value = 581.3685
$ws.Range("N136").Value = -13126.6002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1217.9615
$ws.Range("J5").Value = 1995.5
$ws.Range("L5").Value = 5986.5
$ws.Range("N5").Value = -6210.5

$ws.Range("H113").Value = 750.1698
$ws.Range("I113").Value = 585
$ws.Range("J113").Value = 867.3871
$ws.Range("K113").Value = 1755
$ws.Range("L113").Value = 2602.1613
$ws.Range("M113").Value = 415
$ws.Range("N113").Value = -6942.1613

$ws.Range("H122").Value = 609.15
$ws.Range("I122").Value = 270.5
$ws.Range("J122").Value = 1399.3334
$ws.Range("K122").Value = 2434.5
$ws.Range("L122").Value = 12594.0006
$ws.Range("M122").Value = 15.5
$ws.Range("N122").Value = -17494.0006

$ws.Range("H135").Value = 1217.9615
$ws.Range("J135").Value = 1995.5
$ws.Range("L135").Value = 17959.5
$ws.Range("N135").Value = -23029.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6665.8076
$ws.Range("I70").Value = 7033.75
$ws.Range("J70").Value = 6077.1
$ws.Range("K70").Value = 7033.75
$ws.Range("L70").Value = 6077.1
$ws.Range("M70").Value = -6763.75
$ws.Range("N70").Value = -6617.1

$ws.Range("H73").Value = 6665.8076
$ws.Range("I73").Value = 7033.75
$ws.Range("J73").Value = 6077.1
$ws.Range("K73").Value = 7033.75
$ws.Range("L73").Value = 6077.1
$ws.Range("M73").Value = -6097.75
$ws.Range("N73").Value = -7949.1

$ws.Range("H80").Value = 5234.5884
$ws.Range("I80").Value = 5374.25
$ws.Range("K80").Value = 5374.25
$ws.Range("M80").Value = -4376.25

$ws.Range("H83").Value = 5234.5884
$ws.Range("I83").Value = 5374.25
$ws.Range("K83").Value = 26871.25
$ws.Range("M83").Value = -21879.25

$ws.Range("H102").Value = 2263.2727
$ws.Range("I102").Value = 2110.5881
$ws.Range("K102").Value = 2110.5881
$ws.Range("M102").Value = -488.5880999999999

$ws.Range("H123").Value = 9361.4
$ws.Range("J123").Value = 9361.4
$ws.Range("L123").Value = 9361.4
$ws.Range("N123").Value = -14261.4

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H132").Value = 2096.3242
$ws.Range("I132").Value = 1880.7368
$ws.Range("J132").Value = 2819.1765
$ws.Range("K132").Value = 5642.2104
$ws.Range("L132").Value = 8457.529500000001
$ws.Range("M132").Value = -3112.2104
$ws.Range("N132").Value = -13517.5295

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 16041.857
$ws.Range("I22").Value = 1900.2
$ws.Range("J22").Value = 51396
$ws.Range("K22").Value = 1900.2
$ws.Range("L22").Value = 51396
$ws.Range("M22").Value = -1605.2
$ws.Range("N22").Value = -51986

$ws.Range("H27").Value = 16041.857
$ws.Range("I27").Value = 1900.2
$ws.Range("J27").Value = 51396
$ws.Range("K27").Value = 1900.2
$ws.Range("L27").Value = 51396
$ws.Range("M27").Value = -1793.2
$ws.Range("N27").Value = -51610

$ws.Range("H46").Value = 1652.579
$ws.Range("I46").Value = 1427.1818
$ws.Range("J46").Value = 1962.5
$ws.Range("K46").Value = 1427.1818
$ws.Range("L46").Value = 1962.5
$ws.Range("M46").Value = -1239.1818
$ws.Range("N46").Value = -2338.5

$ws.Range("H55").Value = 250.06667
$ws.Range("I55").Value = 228.75
$ws.Range("J55").Value = 264.27777
$ws.Range("K55").Value = 228.75
$ws.Range("L55").Value = 264.27777
$ws.Range("M55").Value = -55.75
$ws.Range("N55").Value = -610.2777699999999

$ws.Range("H93").Value = 1462.5
$ws.Range("I93").Value = 583.1667
$ws.Range("J93").Value = 2341.8333
$ws.Range("K93").Value = 583.1667
$ws.Range("L93").Value = 2341.8333
$ws.Range("M93").Value = 664.8333
$ws.Range("N93").Value = -4837.8333

$ws.Range("H132").Value = 4163.9287
$ws.Range("I132").Value = 3551.5715
$ws.Range("J132").Value = 4776.2856
$ws.Range("K132").Value = 10654.7145
$ws.Range("L132").Value = 14328.8568
$ws.Range("M132").Value = -8124.7145
$ws.Range("N132").Value = -19388.8568

$ws.Range("H136").Value = 3431.5686
$ws.Range("I136").Value = 2111.6191
$ws.Range("J136").Value = 9591.333000000001
$ws.Range("K136").Value = 6334.8573
$ws.Range("L136").Value = 28773.999
$ws.Range("M136").Value = -3784.8573
$ws.Range("N136").Value = -33873.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10419004
$ws.Range("I132").Value = 17243602
$ws.Range("J132").Value = 2513.5264
$ws.Range("K132").Value = 51730806
$ws.Range("L132").Value = 7540.5792
$ws.Range("M132").Value = -51728276
$ws.Range("N132").Value = -12600.5792
